# correção das notas do fórum para matc65 em 2021.2
# Zera todas as notas/visualizações diárias (colunas B:J) das linhas 3 a 50,
# mantendo a coluna A (matricula) e a linha de cabeçalho (linha 1) e a linha 2 intactas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:J50").Value = 0
